$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.508.34"
$ws.Range("E2").Value = "  +5.57%  "

$ws.Range("D3").Value = "2.292.08"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'303.32"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").Value = "'99.96"
$ws.Range("E6").Value = "  +10.32%  "

$ws.Range("D7").Value = "'0.565"
$ws.Range("E7").Value = "  +1.47%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +5.83%  "

$ws.Range("D10").Value = "'36.28"
$ws.Range("E10").Value = "  +9.67%  "

$ws.Range("D11").Value = "'0.0787"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").Value = "'7.39"
$ws.Range("E12").Value = "  +6.21%  "

$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").Value = "2.643.98"
$ws.Range("E14").Value = "  +3.11%  "

$ws.Range("D15").Value = "2.293.86"
$ws.Range("E15").Value = "  +2.65%  "

$ws.Range("D16").Value = "'13.79"
$ws.Range("E16").Value = "  +3.06%  "

$ws.Range("D17").Value = "'0.811"
$ws.Range("E17").Value = "  +4.48%  "

$ws.Range("D18").Value = "46.521.66"
$ws.Range("E18").Value = "  +6.04%  "

$ws.Range("D19").Value = "'13.01"
$ws.Range("E19").Value = "  +10.40%  "

$ws.Range("D20").Value = "0.0₃0934"
$ws.Range("E20").Value = "  +3.23%  "

$ws.Range("D21").Value = "'5.99"
$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").Value = "'66.04"
$ws.Range("E22").Value = "  +2.99%  "

$ws.Range("D23").Value = "'248.53"
$ws.Range("E23").Value = "  +5.64%  "

$ws.Range("D24").Value = "'2.88"
$ws.Range("E24").Value = "  +2.50%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").Value = "'1.92"
$ws.Range("E26").Value = "  +3.51%  "

$ws.Range("D27").Value = "'42.49"
$ws.Range("E27").Value = "  +7.95%  "

$ws.Range("E28").Value = "  +1.98%  "

$ws.Range("D29").Value = "'9.85"
$ws.Range("E29").Value = "  +5.22%  "

$ws.Range("D30").Value = "'19.96"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").Value = "'2.80"
$ws.Range("E31").Value = "  +12.49%  "

$ws.Range("D32").Value = "'5.62"
$ws.Range("E32").Value = "  +3.69%  "

$ws.Range("D33").Value = "'147.73"
$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("D34").Value = "'0.0791"
$ws.Range("E34").Value = "  +3.55%  "

$ws.Range("D35").Value = "'3.28"
$ws.Range("E35").Value = "  +15.35%  "

$ws.Range("D36").Value = "'0.113"
$ws.Range("E36").Value = "  +9.19%  "

$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  +5.10%  "

$ws.Range("D39").Value = "'15.93"
$ws.Range("E39").Value = "  +19.01%  "

$ws.Range("D40").Value = "'3.96"
$ws.Range("E40").Value = "  +9.95%  "

$ws.Range("D41").Value = "'3.34"
$ws.Range("E41").Value = "  +4.93%  "

$ws.Range("D42").Value = "'0.0301"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "'1.97"
$ws.Range("E44").Value = "  +10.31%  "

$ws.Range("D45").Value = "1.820.98"
$ws.Range("E45").Value = "  +0.98%  "

$ws.Range("D46").Value = "'87.95"
$ws.Range("E46").Value = "  +20.63%  "

$ws.Range("D47").Value = "'0.194"
$ws.Range("E47").Value = "  +5.37%  "

$ws.Range("D48").Value = "'73.16"
$ws.Range("E48").Value = "  +7.69%  "

$ws.Range("D49").Value = "'4.88"
$ws.Range("E49").Value = "  +6.05%  "

$ws.Range("D50").Value = "'95.75"
$ws.Range("E50").Value = "  +1.37%  "

$ws.Range("D51").Value = "2.520.13"
$ws.Range("E51").Value = "  +3.11%  "

